# Script update for Thunderbird LDA Keywords
# Fills in the per-label F1/Precision/Recall/Accuracy score strings for the
# "Count Vectorizer + TFIDF + ngram(3) + POS" configuration on the Ubuntu
# sheet (rows 45-48, columns C:F), then leaves that sheet active/selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ubuntu")

$ws.Range("C45").Value = "0.534 0.420 0.119 0.209 0.728"
$ws.Range("D45").Value = "0.622 0.615 0.667 0.061 0.944 "
$ws.Range("E45").Value = "0.376 0.268 0.063 0.118 0.573"
$ws.Range("F45").Value = "0.789 0.867 0.934 0.852 0.981"

$ws.Range("C46").Value = "0.718 0.727 0.672 0.740 0.751 "
$ws.Range("D46").Value = "0.429 0.291 0.148 0.156 0.623"
$ws.Range("E46").Value = "0.806 0.812 0.614 0.704 0.607"
$ws.Range("F46").Value = "0.687 0.681 0.734 0.775 0.968"

$ws.Range("C47").Value = "0.654 0.574 0.318 0.476 0.784"
$ws.Range("D47").Value = "0.586 0.520 0.271 0.348 0.697"
$ws.Range("E47").Value = "0.521 0.415 0.190 0.316 0.650 "
$ws.Range("F47").Value = "0.791 0.857 0.910 0.930 0.973"

$ws.Range("C48").Value = "0.430 0.273 0.100 0.076 0.721"
$ws.Range("D48").Value = "0.558 0.474 0.476 0.375 0.930"
$ws.Range("E48").Value = "0.279 0.159 0.053 0.039 0.564"
$ws.Range("F48").Value = "0.767 0.850 0.932 0.944 0.980"

# Make Ubuntu the active sheet and move the selection cursor to C49
$ws.Activate()
$ws.Range("C49").Select()
